# Generate Report for Handoff
#
# Updates the localization-status report so the
# d1ba5f85-282e-4c46-a0e2-e03beef41be4 row reflects that its handoff
# package is now ready: Status moves from "In Translation" to
# "Ready for handoff", Priority moves from "ht" (human translation) to
# "mt" (machine translation), and the handoff timestamps are refreshed.
# Column widths on the affected "status" columns are widened so the new,
# longer "Ready for handoff" text fits.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet - row 3 is the d1ba5f85... file
# Columns: A=File Name, B=Path And Name, C=Extension, D=Publish URL,
#          E=zh-cn, F=de-de, G=Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-02 00:22:50"
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# ---------------------------------------------------------------------
# zh-cn sheet - row 3 is the d1ba5f85... file
# Columns: A=Source File Name, B=File Extension, C=Status, D=Source Path,
#          E=Priority, F=Content Duplicate, G=Latest Handoff File,
#          H=Latest Handoff Datetime, ...
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-09-02 00:22:46"
$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797

# ---------------------------------------------------------------------
# de-de sheet - row 3 is the d1ba5f85... file
# Same column layout as zh-cn
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-09-02 00:22:50"
$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
